$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Source Information")

$lines = @(
    "1. The `"Inserted Merged Column`" combines the FirstName, MiddleName, and LastName columns into one column of text values with a space between each name.",
    "2. `"Reordered Columns`" means arranging the columns in the table in a specific order.",
    "3. `"Removed Columns`" means that the columns named `"FirstName`", `"MiddleName`", `"LastName`", `"Gender`" and `"City`" have been taken away from the `"Reordered Columns`" table.",
    "4. Table columns have been renamed from `"Name`" to `"DepartmentName`".",
    "5. This formulated column divides the data in the VacationHours column by 8 to produce a new column with the data for VacationDays.",
    "6. This sentence means that the Custom1 column has been added to the table, and its values are determined by the SickLeaveHours column divided by 8.",
    "7. This sentence adds a new column called `"Round Up`" to the table, and it calculates the rounded up integer of the `"VacationDays`" column and stores it as an integer.",
    "8. Table `"Inserted Round Up`" had its column `"Round Up`" removed, resulting in `"Removed Columns1`".",
    "9. Rounded Off means to round the number of Vacation Days to two decimal places.",
    "10. The code rounds the `"Sick Leave Days`" column to 2 decimal places and changes the data type to a number.",
    "11. Table.RemoveColumns removes the columns `"VacationHours`" and `"SickLeaveHours`" from the table `"RoundedOff1`".",
    "12. Reorder Columns1 rearranges the Removed Columns2 table so that the columns are in the specified order: EmployeeName, DepartmentName, VacationDays, SickLeaveDays, and MaritalStatus.",
    "13. This sentence means that the MaritalStatus column has been removed from the Reordered Columns1 table."
)

$text = ($lines -join "`n`n") + "`n`n"
$ws.Range("I2").Value = $text
